# Actualización automática 2025-09-22 08:22:24
$wb = $excel.ActiveWorkbook

$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# VENTAS POR GRUPO: commission / sales value for RIOS CARRION ANGEL BENIGNO / TOSCANO RAMIREZ MONICA CECILIA
$wsVentasGrupo.Range("M22").Value = -347.92

# VENTA MENSUAL: same value repeated, plus recomputed total row
$wsVentaMensual.Range("F22").Value = -347.92
$wsVentaMensual.Range("F26").Value = 19557.03

# CUMPLIMIENTO MENSUAL: PORCELANATO row (12) and TOTAL row (15) recomputed
$wsCumplimiento.Range("D12").Value = 19447.67
$wsCumplimiento.Range("E12").Value = 23652.4154117774
$wsCumplimiento.Range("F12").Value = 0.4512211475730808

$wsCumplimiento.Range("D15").Value = 19557.03
$wsCumplimiento.Range("E15").Value = 38646.43623249458
$wsCumplimiento.Range("F15").Value = 0.3360114313790035
